# Implement "Incoming transfer during period" breakdown row (row 13) on the
# "File active" worksheet, mirroring the per-gender/age-band structure used
# by the other indicator rows (rows 12, 18, 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("File active")
$ws.Activate()

$ws.Range("B13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', gender: 0}'
$ws.Range("C13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', gender: 1}'
$ws.Range("D13").Value = '{key:''INCOMING_TRANSFER_DURING_PERIOD'', age_max: 15}'
$ws.Range("E13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 15}'
$ws.Range("F13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD''}'
$ws.Range("G13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_max: 15, gender: 0}'
$ws.Range("H13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_max: 15, gender: 1}'
$ws.Range("I13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 15, gender: 0}'
$ws.Range("J13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 15, gender: 1}'
$ws.Range("K13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_max: 1, gender: 0}'
$ws.Range("L13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_max: 1, gender: 1}'
$ws.Range("M13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 1, age_max: 4, gender: 0}'
$ws.Range("N13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_max: 1, gender: 0}'
$ws.Range("O13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 5, age_max: 9, gender: 0}'
$ws.Range("P13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 5, age_max: 9, gender: 1}'
$ws.Range("Q13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 10, age_max: 14, gender: 0}'
$ws.Range("R13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 10, age_max: 14, gender: 1}'
$ws.Range("S13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 15, age_max: 19, gender: 0}'
$ws.Range("T13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 15, age_max: 19, gender: 1}'
$ws.Range("U13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 20, age_max: 24, gender: 0}'
$ws.Range("V13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 20, age_max: 24, gender: 1}'
$ws.Range("W13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 25, age_max: 49, gender: 0}'
$ws.Range("X13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 25, age_max: 49, gender: 1}'
$ws.Range("Y13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 50, gender: 0}'
$ws.Range("Z13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', age_min: 50, gender: 1}'
$ws.Range("AA13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', gender: 0}'
$ws.Range("AB13").Value = '{key: ''INCOMING_TRANSFER_DURING_PERIOD'', gender: 1}'

# The row grows taller now that it carries wrapped header text in every
# column (matches the height used for the analogous rows 18/20).
$ws.Rows.Item(13).RowHeight = 29.85

# Restore the editor's scroll position / selection as left by the author.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 9
$ws.Range("AC13").Select()
